{"js": "// Add six new bulleted entries ([Musai], [Tur], [Altell], [Altair],\n// [Kalabeth], [Endor]) to the \"aynu code atoms\" list, right after the\n// existing \"[Valaktioth]\" entry.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph (\"[Valaktioth]\") that the new items follow.\nconst anchor = paragraphs.items.find((p) => p.text === \"[Valaktioth]\");\nif (!anchor) {\n  throw new Error('Could not find the \"[Valaktioth]\" paragraph to anchor the new list items.');\n}\n\nconst newItems = [\"[Musai]\", \"[Tur]\", \"[Altell]\", \"[Altair]\", \"[Kalabeth]\", \"[Endor]\"];\n\n// Insert each new bullet item right after the anchor, in order, chaining off\n// the previously inserted paragraph so the final order matches newItems.\nlet previous = anchor;\nfor (const text of newItems) {\n  previous = previous.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Add six new bulleted entries ([Musai], [Tur], [Altell], [Altair],\n# [Kalabeth], [Endor]) to the \"aynu code atoms\" list, right after the\n# existing \"[Valaktioth]\" entry.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"[Valaktioth]\") that the new items follow.\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"[Valaktioth]`r\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw 'Could not find the \"[Valaktioth]\" paragraph to anchor the new list items.'\n}\n\n$newItems = @(\"[Musai]\", \"[Tur]\", \"[Altell]\", \"[Altair]\", \"[Kalabeth]\", \"[Endor]\")\n\n# Insert each new bullet item right after the anchor, in order, chaining off\n# the previously inserted paragraph so the final order matches $newItems.\n$current = $anchor\nforeach ($text in $newItems) {\n    $current.Range.InsertParagraphAfter()\n    $current = $current.Next()\n    $current.Range.Text = $text\n}\n"}
